# Calculator Keyboard Layout - apply "Add bitmaps for numbers" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet view: scroll back to top-left (drop topLeftCell="A4") and move the
# selection from L3 to U9. Selecting the target cell both moves the active
# selection and resets the window's topLeftCell to the default.
$ws.Range("U9").Select()

# --- S5 used to read "Frac"; it now shows a stacked "a/b" glyph (superscript
# "a", plain "/", subscript "b") matching the style of the other stacked
# fraction-like glyphs already in the workbook (e.g. the "a b/c" key).
$s5 = $ws.Range("S5")
$s5.Value = "a/b"

$s5a = $s5.Characters(1, 1)
$s5a.Font.Name = "Consolas"
$s5a.Font.Size = 14
$s5a.Font.Superscript = $true

$s5slash = $s5.Characters(2, 1)
$s5slash.Font.Name = "Consolas"
$s5slash.Font.Size = 14

$s5b = $s5.Characters(3, 1)
$s5b.Font.Name = "Consolas"
$s5b.Font.Size = 14
$s5b.Font.Subscript = $true

# --- AB6 used to read "EE"; it now matches F6 exactly (same rich "x10^x"
# glyph, same style s=19) - copy F6 (value + formatting) onto AB6 so the
# cell ends up sharing F6's formatted string + style.
$ws.Range("F6").Copy($ws.Range("AB6"))

Write-Host "edit applied"
